$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value2 = 7137.885
$ws.Range("I106").Value2 = 2060.8462
$ws.Range("K106").Value2 = 2060.8462
$ws.Range("M106").Value2 = -1429.8462
$ws.Range("H112").Value2 = 2205.8823
$ws.Range("I112").Value2 = 1000
$ws.Range("J112").Value2 = 3050
$ws.Range("K112").Value2 = 3000
$ws.Range("L112").Value2 = 9150
$ws.Range("M112").Value2 = -1892
$ws.Range("N112").Value2 = -11366
$ws.Range("H113").Value2 = 3010.6667
$ws.Range("I113").Value2 = 3065
$ws.Range("J113").Value2 = 2902
$ws.Range("K113").Value2 = 3065
$ws.Range("L113").Value2 = 2902
$ws.Range("M113").Value2 = 189
$ws.Range("N113").Value2 = -9410
$ws.Range("H116").Value2 = 4710.5654
$ws.Range("I116").Value2 = 4108.353
$ws.Range("J116").Value2 = 6416.8335
$ws.Range("K116").Value2 = 4108.353
$ws.Range("L116").Value2 = 6416.8335
$ws.Range("M116").Value2 = -666.3530000000001
$ws.Range("N116").Value2 = -13300.8335
$ws.Range("H137").Value2 = 8871
$ws.Range("I137").Value2 = 4099.8623
$ws.Range("J137").Value2 = 13334.322
$ws.Range("K137").Value2 = 12299.5869
$ws.Range("L137").Value2 = 40002.966
$ws.Range("M137").Value2 = -9749.586899999998
$ws.Range("N137").Value2 = -45102.966
$ws.Range("H138").Value2 = 5921.8945
$ws.Range("J138").Value2 = 5742.2856
$ws.Range("L138").Value2 = 17226.8568
$ws.Range("N138").Value2 = -27506.8568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 9447599
$ws.Range("I32").Value2 = 11120587
$ws.Range("K32").Value2 = 11120587
$ws.Range("M32").Value2 = -11120300
$ws.Range("H61").Value2 = 13920175
$ws.Range("I61").Value2 = 19233124
$ws.Range("J61").Value2 = 106505.1
$ws.Range("K61").Value2 = 19233124
$ws.Range("L61").Value2 = 106505.1
$ws.Range("M61").Value2 = -19232912
$ws.Range("N61").Value2 = -106929.1
$ws.Range("H110").Value2 = 1442.8334
$ws.Range("I110").Value2 = 1442.8334
$ws.Range("K110").Value2 = 1442.8334
$ws.Range("M110").Value2 = 602.1666
$ws.Range("H136").Value2 = 13920175
$ws.Range("I136").Value2 = 19233124
$ws.Range("J136").Value2 = 106505.1
$ws.Range("K136").Value2 = 57699372
$ws.Range("L136").Value2 = 319515.3
$ws.Range("M136").Value2 = -57696822
$ws.Range("N136").Value2 = -324615.3

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 26674.768
$ws.Range("I134").Value2 = 3264.5334
$ws.Range("J134").Value2 = 80698.38
$ws.Range("K134").Value2 = 9793.600199999999
$ws.Range("L134").Value2 = 242095.14
$ws.Range("M134").Value2 = -7258.600199999999
$ws.Range("N134").Value2 = -247165.14

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 16498.75
$ws.Range("I16").Value2 = 21331.666
$ws.Range("J16").Value2 = 2000
$ws.Range("K16").Value2 = 21331.666
$ws.Range("L16").Value2 = 2000
$ws.Range("M16").Value2 = -21044.666
$ws.Range("N16").Value2 = -2574
$ws.Range("H31").Value2 = 348002.94
$ws.Range("I31").Value2 = 10016.3
$ws.Range("J31").Value2 = 488830.7
$ws.Range("K31").Value2 = 10016.3
$ws.Range("L31").Value2 = 488830.7
$ws.Range("M31").Value2 = -9721.299999999999
$ws.Range("N31").Value2 = -489420.7
$ws.Range("H34").Value2 = 348002.94
$ws.Range("I34").Value2 = 10016.3
$ws.Range("J34").Value2 = 488830.7
$ws.Range("K34").Value2 = 10016.3
$ws.Range("L34").Value2 = 488830.7
$ws.Range("M34").Value2 = -9814.299999999999
$ws.Range("N34").Value2 = -489234.7
$ws.Range("H113").Value2 = 16498.75
$ws.Range("I113").Value2 = 21331.666
$ws.Range("J113").Value2 = 2000
$ws.Range("K113").Value2 = 21331.666
$ws.Range("L113").Value2 = 2000
$ws.Range("M113").Value2 = -19161.666
$ws.Range("N113").Value2 = -6340
$ws.Range("H134").Value2 = 401426.6
$ws.Range("I134").Value2 = 527640.4399999999
$ws.Range("J134").Value2 = 1749.5
$ws.Range("K134").Value2 = 1582921.32
$ws.Range("L134").Value2 = 5248.5
$ws.Range("M134").Value2 = -1580386.32
$ws.Range("N134").Value2 = -10318.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value2 = 913.4737
$ws.Range("I60").Value2 = 449
$ws.Range("J60").Value2 = 1079.3572
$ws.Range("K60").Value2 = 1347
$ws.Range("L60").Value2 = 3238.0716
$ws.Range("M60").Value2 = -1096
$ws.Range("N60").Value2 = -3740.0716
$ws.Range("H68").Value2 = 3422.9023
$ws.Range("I68").Value2 = 3171.818
$ws.Range("J68").Value2 = 3514.9666
$ws.Range("K68").Value2 = 9515.454000000002
$ws.Range("L68").Value2 = 10544.8998
$ws.Range("M68").Value2 = -8704.454000000002
$ws.Range("N68").Value2 = -12166.8998
$ws.Range("H71").Value2 = 3422.9023
$ws.Range("I71").Value2 = 3171.818
$ws.Range("J71").Value2 = 3514.9666
$ws.Range("K71").Value2 = 28546.362
$ws.Range("L71").Value2 = 31634.6994
$ws.Range("M71").Value2 = -24490.362
$ws.Range("N71").Value2 = -39746.6994

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value2 = 3710.6155
$ws.Range("I113").Value2 = 2889.5715
$ws.Range("K113").Value2 = 2889.5715
$ws.Range("M113").Value2 = -719.5715
$ws.Range("H122").Value2 = 2506.3076
$ws.Range("I122").Value2 = 1759.96
$ws.Range("J122").Value2 = 3839.0715
$ws.Range("K122").Value2 = 5279.88
$ws.Range("L122").Value2 = 11517.2145
$ws.Range("M122").Value2 = -2829.88
$ws.Range("N122").Value2 = -16417.2145
$ws.Range("H132").Value2 = 37044772
$ws.Range("I132").Value2 = 45459108
$ws.Range("J132").Value2 = 21699.6
$ws.Range("K132").Value2 = 136377324
$ws.Range("L132").Value2 = 65098.8
$ws.Range("M132").Value2 = -136374794
$ws.Range("N132").Value2 = -70158.79999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 678.3333
$ws.Range("I22").Value2 = 773.3333
$ws.Range("J22").Value2 = 583.3333
$ws.Range("K22").Value2 = 773.3333
$ws.Range("L22").Value2 = 583.3333
$ws.Range("M22").Value2 = -478.3333
$ws.Range("N22").Value2 = -1173.3333
$ws.Range("H27").Value2 = 678.3333
$ws.Range("I27").Value2 = 773.3333
$ws.Range("J27").Value2 = 583.3333
$ws.Range("K27").Value2 = 773.3333
$ws.Range("L27").Value2 = 583.3333
$ws.Range("M27").Value2 = -666.3333
$ws.Range("N27").Value2 = -797.3333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 16673375
$ws.Range("I62").Value2 = 6001
$ws.Range("K62").Value2 = 6001
$ws.Range("M62").Value2 = -5377
$ws.Range("H65").Value2 = 16673375
$ws.Range("I65").Value2 = 6001
$ws.Range("K65").Value2 = 30005
$ws.Range("M65").Value2 = -26885
$ws.Range("H81").Value2 = 139999.67
$ws.Range("J81").Value2 = 110000
$ws.Range("L81").Value2 = 220000
$ws.Range("N81").Value2 = -222122
$ws.Range("H84").Value2 = 139999.67
$ws.Range("J84").Value2 = 110000
$ws.Range("L84").Value2 = 1100000
$ws.Range("N84").Value2 = -1110608
